$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3183.25
$ws.Range("I76").Value = 3028.4285
$ws.Range("K76").Value = 3028.4285
$ws.Range("M76").Value = -2713.4285
$ws.Range("H79").Value = 3183.25
$ws.Range("I79").Value = 3028.4285
$ws.Range("K79").Value = 3028.4285
$ws.Range("M79").Value = -1936.4285
$ws.Range("H106").Value = 9010803
$ws.Range("I106").Value = 12346804
$ws.Range("K106").Value = 12346804
$ws.Range("M106").Value = -12346173
$ws.Range("H113").Value = 4700.24
$ws.Range("I113").Value = 4388.8887
$ws.Range("K113").Value = 4388.8887
$ws.Range("M113").Value = -1134.8887
$ws.Range("H129").Value = 296256.62
$ws.Range("I129").Value = 394
$ws.Range("J129").Value = 359655.75
$ws.Range("K129").Value = 1182
$ws.Range("L129").Value = 1078967.25
$ws.Range("M129").Value = 3818
$ws.Range("N129").Value = -1088967.25
$ws.Range("H137").Value = 47665.91
$ws.Range("I137").Value = 2079.9
$ws.Range("J137").Value = 85654.25
$ws.Range("K137").Value = 6239.700000000001
$ws.Range("L137").Value = 256962.75
$ws.Range("M137").Value = -3689.700000000001
$ws.Range("N137").Value = -262062.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15099.223
$ws.Range("I32").Value = 16313.9375
$ws.Range("J32").Value = 5381.5
$ws.Range("K32").Value = 16313.9375
$ws.Range("L32").Value = 5381.5
$ws.Range("M32").Value = -16026.9375
$ws.Range("N32").Value = -5955.5
$ws.Range("H110").Value = 407.875
$ws.Range("I110").Value = 415
$ws.Range("K110").Value = 415
$ws.Range("M110").Value = 1630
$ws.Range("H122").Value = 1811.9459
$ws.Range("I122").Value = 1420.3462
$ws.Range("K122").Value = 4261.0386
$ws.Range("M122").Value = -1811.0386
$ws.Range("H132").Value = 47654.184
$ws.Range("I132").Value = 2095.1924
$ws.Range("K132").Value = 6285.5772
$ws.Range("M132").Value = -3755.5772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 16641.031
$ws.Range("I58").Value = 983.76666
$ws.Range("K58").Value = 983.76666
$ws.Range("M58").Value = -780.76666
$ws.Range("H136").Value = 16641.031
$ws.Range("I136").Value = 983.76666
$ws.Range("K136").Value = 2951.29998
$ws.Range("M136").Value = -401.2999799999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3545.8
$ws.Range("J114").Value = 3106
$ws.Range("L114").Value = 9318
$ws.Range("N114").Value = -15826
$ws.Range("H118").Value = 35718704
$ws.Range("J118").Value = 8556.714
$ws.Range("L118").Value = 25670.142
$ws.Range("N118").Value = -28156.142
$ws.Range("H122").Value = 717.3
$ws.Range("I122").Value = 367.5
$ws.Range("J122").Value = 804.75
$ws.Range("K122").Value = 3307.5
$ws.Range("L122").Value = 7242.75
$ws.Range("M122").Value = -857.5
$ws.Range("N122").Value = -12142.75
$ws.Range("H131").Value = 766.2
$ws.Range("J131").Value = 797.55316
$ws.Range("L131").Value = 2392.65948
$ws.Range("N131").Value = -12472.65948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4130.769
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 4190.909
$ws.Range("K80").Value = 3800
$ws.Range("L80").Value = 4190.909
$ws.Range("M80").Value = -2802
$ws.Range("N80").Value = -6186.909
$ws.Range("H83").Value = 4130.769
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 4190.909
$ws.Range("K83").Value = 19000
$ws.Range("L83").Value = 20954.545
$ws.Range("M83").Value = -14008
$ws.Range("N83").Value = -30938.545
$ws.Range("H102").Value = 35718024
$ws.Range("I102").Value = 45458620
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 45458620
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -45456998
$ws.Range("N102").Value = -5744
$ws.Range("H126").Value = 4311.364
$ws.Range("I126").Value = 3630.55
$ws.Range("K126").Value = 10891.65
$ws.Range("M126").Value = -8421.650000000001
$ws.Range("H132").Value = 66361.336
$ws.Range("I132").Value = 56830.21
$ws.Range("K132").Value = 170490.63
$ws.Range("M132").Value = -167960.63

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1718.8334
$ws.Range("I22").Value = 1414.7
$ws.Range("J22").Value = 2099
$ws.Range("K22").Value = 1414.7
$ws.Range("L22").Value = 2099
$ws.Range("M22").Value = -1119.7
$ws.Range("N22").Value = -2689
$ws.Range("H27").Value = 1718.8334
$ws.Range("I27").Value = 1414.7
$ws.Range("J27").Value = 2099
$ws.Range("K27").Value = 1414.7
$ws.Range("L27").Value = 2099
$ws.Range("M27").Value = -1307.7
$ws.Range("N27").Value = -2313
$ws.Range("H46").Value = 969.5
$ws.Range("J46").Value = 939.8
$ws.Range("L46").Value = 939.8
$ws.Range("N46").Value = -1315.8
$ws.Range("H55").Value = 301.33334
$ws.Range("I55").Value = 327.5
$ws.Range("J55").Value = 280.4
$ws.Range("K55").Value = 327.5
$ws.Range("L55").Value = 280.4
$ws.Range("M55").Value = -154.5
$ws.Range("N55").Value = -626.4
$ws.Range("H61").Value = 6137
$ws.Range("I61").Value = 2917.8333
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 2917.8333
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -2715.8333
$ws.Range("N61").Value = -10404
$ws.Range("H113").Value = 6137
$ws.Range("I113").Value = 2917.8333
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2917.8333
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -747.8332999999998
$ws.Range("N113").Value = -14340
$ws.Range("H122").Value = 936373.4399999999
$ws.Range("I122").Value = 1963214.2
$ws.Range("J122").Value = 2881.7273
$ws.Range("K122").Value = 5889642.6
$ws.Range("L122").Value = 8645.1819
$ws.Range("M122").Value = -5887192.6
$ws.Range("N122").Value = -13545.1819
$ws.Range("H132").Value = 1558.25
$ws.Range("I132").Value = 1132.3871
$ws.Range("K132").Value = 3397.1613
$ws.Range("M132").Value = -867.1612999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 50000748
$ws.Range("I81").Value = 818
$ws.Range("J81").Value = 166667250
$ws.Range("K81").Value = 1636
$ws.Range("L81").Value = 333334500
$ws.Range("M81").Value = -575
$ws.Range("N81").Value = -333336622
$ws.Range("H84").Value = 50000748
$ws.Range("I84").Value = 818
$ws.Range("J84").Value = 166667250
$ws.Range("K84").Value = 8180
$ws.Range("L84").Value = 1666672500
$ws.Range("M84").Value = -2876
$ws.Range("N84").Value = -1666683108
